# CryCompanywiseStockReport_1 — stock-count correction pass.
#
# A batch of "Closing Qty" (column F) / "Closing Value" (column G = F * D,
# Cost Price) corrections for individual SKUs, the resulting "Sub Total:"
# (column B) roll-ups for the affected company blocks, and a couple of
# adjacent data-row reorderings (same two SKUs, values swapped between the
# two rows) deeper in the sheet. Grand-total rows at the bottom are updated
# to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Simple quantity/value corrections (Qty -> lower, Value = Qty * Cost) ---

$ws.Range("F30").Value = 22
$ws.Range("G30").Value = 2253.9

$ws.Range("F118").Value = 34
$ws.Range("G118").Value = 2685.32

$ws.Range("F119").Value = 104
$ws.Range("G119").Value = 14596.4

$ws.Range("F126").Value = 78
$ws.Range("G126").Value = 10512.06

$ws.Range("F155").Value = 49
$ws.Range("G155").Value = 8883.209999999999

$ws.Range("F160").Value = 289
$ws.Range("G160").Value = 9635.26

$ws.Range("F178").Value = 79
$ws.Range("G178").Value = 5018.08

$ws.Range("F258").Value = 64
$ws.Range("G258").Value = 5705.6

$ws.Range("F286").Value = 8
$ws.Range("G286").Value = 776.8

$ws.Range("F439").Value = 87
$ws.Range("G439").Value = 838.6799999999999

$ws.Range("F536").Value = 13
$ws.Range("G536").Value = 561.34

$ws.Range("F557").Value = 5
$ws.Range("G557").Value = 3725.95

$ws.Range("F631").Value = 295
$ws.Range("G631").Value = 10864.85

$ws.Range("F680").Value = 402
$ws.Range("G680").Value = 65570.22

$ws.Range("F688").Value = 0
$ws.Range("G688").Value = 0

# --- Adjacent-row reorderings (Code/Item/Price/Qty/Value swapped between the two rows) ---

# Rows 167 <-> 168 (COL-Colgate Zigzag Charcoal Pack of 4 Toothbrush)
$ws.Range("B167").Value = 64350
$ws.Range("E167").Value = 70.63
$ws.Range("F167").Value = 2
$ws.Range("G167").Value = 132.88

$ws.Range("B168").Value = 57756
$ws.Range("E168").Value = 79.37
$ws.Range("F168").Value = -100
$ws.Range("G168").Value = -6644

# Rows 298 <-> 299 (HIM-Total Care Baby Pants Drapers-Xl-9S, two case variants)
$ws.Range("B298").Value = 64985
$ws.Range("C298").Value = "HIM-TOTAL CARE BABY PANTS DRAPERS-XL-9S"
$ws.Range("F298").Value = 12
$ws.Range("G298").Value = 1052.4

$ws.Range("B299").Value = 66196
$ws.Range("C299").Value = "HIM-Total Care Baby Pants Drapers-Xl-9S"
$ws.Range("F299").Value = 1
$ws.Range("G299").Value = 87.7

# Rows 303 <-> 304
$ws.Range("B303").Value = 61610
$ws.Range("E303").Value = 122.71
$ws.Range("F303").Value = -58
$ws.Range("G303").Value = -5957.18

$ws.Range("B304").Value = 63565
$ws.Range("E304").Value = 109.19
$ws.Range("F304").Value = 60
$ws.Range("G304").Value = 6162.6

# Rows 312 <-> 313
$ws.Range("B312").Value = 63531
$ws.Range("E312").Value = 152.53
$ws.Range("F312").Value = 24
$ws.Range("G312").Value = 3443.52

$ws.Range("B313").Value = 57802
$ws.Range("E313").Value = 162.71
$ws.Range("F313").Value = -79
$ws.Range("G313").Value = -11334.92

# Rows 387 <-> 388
$ws.Range("B387").Value = 47097
$ws.Range("D387").Value = 112.28
$ws.Range("E387").Value = 134.16
$ws.Range("F387").Value = 15
$ws.Range("G387").Value = 1684.2

$ws.Range("B388").Value = 58047
$ws.Range("D388").Value = 105.54
$ws.Range("E388").Value = 126.1
$ws.Range("F388").Value = 32
$ws.Range("G388").Value = 3377.28

# Rows 502 <-> 503
$ws.Range("B502").Value = 64833
$ws.Range("E502").Value = 34.9
$ws.Range("F502").Value = 88
$ws.Range("G502").Value = 2889.04

$ws.Range("B503").Value = 60025
$ws.Range("E503").Value = 37.22
$ws.Range("F503").Value = -98
$ws.Range("G503").Value = -3217.34

# Rows 512 <-> 513
$ws.Range("B512").Value = 60022
$ws.Range("E512").Value = 37.22
$ws.Range("F512").Value = -113
$ws.Range("G512").Value = -3709.79

$ws.Range("B513").Value = 64830
$ws.Range("E513").Value = 34.9
$ws.Range("F513").Value = 83
$ws.Range("G513").Value = 2724.89

# --- "Sub Total:" roll-ups for every company block touched above ---

$ws.Range("B40").Value = 53252.68
$ws.Range("B129").Value = 67631.19
$ws.Range("B161").Value = 33855.87
$ws.Range("B199").Value = 55643.54
$ws.Range("B301").Value = 95811.39999999999
$ws.Range("B444").Value = 20827.04
$ws.Range("B541").Value = 19468.47
$ws.Range("B562").Value = 35699.85
$ws.Range("B634").Value = 192659.23
$ws.Range("B686").Value = 66582.77
$ws.Range("B690").Value = 17390.48

# --- Sheet-wide "Sub Total:" / "Grand Total:" at the bottom ---

$ws.Range("B724").Value = 2264682.91
$ws.Range("B725").Value = 2264682.91
